# Update country data files
# Adds a new "MSME definitions" table (rows 23-27) to the Greece Summary
# sheet, and re-adds the source citation two rows lower (rows 32-33),
# duplicating the shared strings exactly as the target workbook does.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New table header (row 23) ------------------------------------------
$ws.Range("B23").Value = "Number of employees"
$ws.Range("B23").Style = "title"
$ws.Range("C23").Value = "Assets (local currency, unless noted otherwise)"
$ws.Range("C23").Style = "title"
$ws.Range("D23").Value = "Turnover (local currency, unless noted otherwise)"
$ws.Range("D23").Style = "title"

# --- Micro row (24) -------------------------------------------------------
$ws.Range("A24").Value = "Micro"
$ws.Range("A24").Style = "Normal"
$ws.Range("B24").Value = "<10"
$ws.Range("B24").Style = "Normal"
$ws.Range("C24").Value = "≤ €2 Millionlion (previously not defined)"
$ws.Range("C24").Style = "Normal"
$ws.Range("D24").Value = "≤ € 2 Millionlion"
$ws.Range("D24").Style = "Normal"

# --- Small row (25) ---------------------------------------------------------
$ws.Range("A25").Value = "Small"
$ws.Range("A25").Style = "Normal"
$ws.Range("B25").Value = "<50"
$ws.Range("B25").Style = "Normal"
$ws.Range("C25").Value = "≤ €10 Millionlion (in 1996 €5 Millionlion)"
$ws.Range("C25").Style = "Normal"
$ws.Range("D25").Value = "≤ € 10 Millionlion"
$ws.Range("D25").Style = "Normal"

# --- Medium row (26) --------------------------------------------------------
$ws.Range("A26").Value = "Medium"
$ws.Range("A26").Style = "Normal"
$ws.Range("B26").Value = "<250"
$ws.Range("B26").Style = "Normal"
$ws.Range("C26").Value = "≤ €43 Millionlion (in 1996 € 27 Millionlion)"
$ws.Range("C26").Style = "Normal"
$ws.Range("D26").Value = "≤ € 50 Millionlion"
$ws.Range("D26").Style = "Normal"

# --- Large row (27) ----------------------------------------------------------
$ws.Range("A27").Value = "Large"
$ws.Range("A27").Style = "Normal"
$ws.Range("B27").Value = ">249"
$ws.Range("B27").Style = "Normal"
$ws.Range("C27").Value = "> €43 Millionlion (in 1996 € 27 Millionlion)"
$ws.Range("C27").Style = "Normal"
$ws.Range("D27").Value = "> € 50 Millionlion"
$ws.Range("D27").Style = "Normal"

# --- Source citation, now two rows further down (32-33) ----------------------
$ws.Range("A32").Value = "SME Performance Review EU"
$ws.Range("A32").Style = "title"
$ws.Range("A33").Value = 'SME Performance Review EU, "SBA Fact sheet", 2013.  Available at http://ec.europa.eu/enterprise/policies/sme/facts-figures-analysis/performance-review/index_en.htm'
$ws.Range("A33").Style = "source"

$ws.Range("A1:D33").Value | Out-Null
